$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet2 (DPLKKPS143-002): insert a new "NO_REGISTER" column before column N ---
$ws2.Columns("N:N").Insert()

# New column header + value
$ws2.Range("N1").Value = "NO_REGISTER"
$ws2.Range("N2").Value = "M11220800000027"

# Clear the (now shifted) NO_PESERTA value in P2 - this scenario no longer records it
$ws2.Range("P2").ClearContents()

# Update the (now shifted) SNF values for the new scenario
$ws2.Range("Q2").Value = 308000
$ws2.Range("R2").Value = 1292000

# Update preparation text (long scripted note) and USERID for the new run
$ws2.Range("F2").Value = "Username : 33028;" + [char]10 + "Password : bni1234;" + [char]10 + "Role : 10 - Asisten Settlement;" + [char]10 + "Keterangan Perubahan : KEP.TRX.445 melakukan Split Iuran;" + [char]10 + "Saldo Nominal Final - Saldo Awal Iuran Pribadi : 308.000,00;" + [char]10 + "Saldo Nominal Final - Saldo Awal Iuran Perusahaan : 1.292.000,00;" + [char]10 + "Saldo Nominal Final - Saldo Awal Iuran Sukarela : 0,00;" + [char]10 + "Saldo Nominal Final - Saldo Awal Pengalihan Iuran Karyawan : 0,00;" + [char]10 + "Saldo Nominal Final - Saldo Awal Pengalihan Iuran Perusahaan : 0,00;" + [char]10 + "Status Register : 1 - Lanjutkan Ke Verifikasi;" + [char]10 + "Keterangan Register : KEP.TRX.445 Lanjutkan Verifikasi"
$ws2.Range("G2").Value = 33028

# Row 2 grows taller to fit the longer wrapped preparation text
$ws2.Rows(2).RowHeight = 242.25

# --- Sheet1 (DPLKKPS143-001): selection moves from S2 to G2 ---
$ws1.Range("G2").Select()

# --- Sheet2 becomes the active/selected tab, selection moves to R2 ---
$ws2.Activate()
$ws2.Range("R2").Select()
